$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading-percent values for columns B,C,D,E,G,H,I,J,L across rows 2-25
# (columns F,K,M,N,O are unchanged zeros). Row order below follows A2..A25 (index 0..23).
$newValues = @(
    @(23.35199135438161, 9.040056169400552, 12.972983226917, 12.68481465606163, 60.65750051069526, 22.10401503027334, 33.79039134410352, 7.84161638198797, 13.8230676791848),
    @(23.00377652425601, 8.630193828868856, 12.98364283608631, 12.71412108419088, 60.57672846018547, 22.14630071114788, 33.87820201712407, 7.842492910190352, 13.81981706949494),
    @(22.79273218477701, 8.370456018728326, 12.99247626283807, 12.73342881805583, 60.54551340334883, 22.17695490009561, 33.93970114945091, 7.843117112492356, 13.82019276451908),
    @(22.70753268043791, 8.262741364766462, 12.99665105588811, 12.74162761169938, 60.53740803003073, 22.1906222073212, 33.96666113146141, 7.843393142134711, 13.82094313223591),
    @(22.6934370066258, 8.244747709947466, 12.99737899889913, 12.74300900411304, 60.53634059107619, 22.19296252694115, 33.97125223991668, 7.843440286001155, 13.82110382560138),
    @(22.7915797614176, 8.36901067914356, 12.99253023755553, 12.73353805011742, 60.54538541783936, 22.17713446831976, 33.94005706512272, 7.843120747368597, 13.82020046491018),
    @(23.23141914524869, 8.90050358360458, 12.97618354930515, 12.69464717430126, 60.62583307430958, 22.11761925150011, 33.81908974029925, 7.841900777523149, 13.82145526620912),
    @(24.1110662274246, 9.872570689475427, 12.96229506420098, 12.62878385815318, 60.92951281340652, 22.03829873016448, 33.64239858976044, 7.840189502626811, 13.84268458582765),
    @(24.76141858472775, 10.53716848423547, 12.96317300074061, 12.58670753617236, 61.24141292362361, 22.00303536151123, 33.54995488491696, 7.839345779737174, 13.86964556922238),
    @(25.0568498374114, 10.82764625663934, 12.9659767029638, 12.56893106761256, 61.40244728876924, 21.99203187898373, 33.5161095692705, 7.839051421864593, 13.88435642626909),
    @(25.16855584741199, 10.93586557069273, 12.9673836036027, 12.56239533035292, 61.46616099726905, 21.98859247606472, 33.50448068476506, 7.838952792098647, 13.8902764430232),
    @(25.14450707498177, 10.91263874951557, 12.96706526011915, 12.56379421448713, 61.45231790899146, 21.98930081965177, 33.50693224158807, 7.838973463260455, 13.88898596618191),
    @(25.06604391495308, 10.83658554677482, 12.9660855349682, 12.56838944594364, 61.40763440311503, 21.99173432394611, 33.51512901203665, 7.839043050407003, 13.88483648325925),
    @(25.01795803945364, 10.78976718079192, 12.96553036098597, 12.57122964976836, 61.38061976685879, 21.9933197226898, 33.52030464235003, 7.839087345550663, 13.88234022232511),
    @(24.74209379909297, 10.51793998634376, 12.96303811150715, 12.58789669227827, 61.23127266860627, 22.00385609922492, 33.55233255321193, 7.839366814516627, 13.86873320247335),
    @(24.57267562151522, 10.34809075511685, 12.96212467767666, 12.59847054975362, 61.14454484948932, 22.01161235552426, 33.57408819006714, 7.839561152387355, 13.86101062232405),
    @(24.47519770499689, 10.24928608033976, 12.96182564523414, 12.60468078446631, 61.09646525177942, 22.01654759506003, 33.58737359879218, 7.839681351939245, 13.85679920586389),
    @(24.44219100908409, 10.21564400859162, 12.9617632872758, 12.60680553232724, 61.08049669282566, 22.01829990599465, 33.59200419614757, 7.839723496670691, 13.85541293463792),
    @(24.59071468694624, 10.36628714770081, 12.96219849118205, 12.59733165499572, 61.15359060350498, 22.01073760489397, 33.57169229760847, 7.839539593384594, 13.86180887459651),
    @(25.08909582765654, 10.85897301022887, 12.96636394060524, 12.56703440321492, 61.42068503006289, 21.99099978226723, 33.51268913419002, 7.839022262785396, 13.88604582678093),
    @(25.4137938960939, 11.17058009418917, 12.97109800688126, 12.54837465897376, 61.61116687397903, 21.98234080981632, 33.48105253048302, 7.838758970477196, 13.90392095144125),
    @(25.24062451983936, 11.00524182050156, 12.96838750204374, 12.55822940553449, 61.50805431179488, 21.98657334875956, 33.49730158871058, 7.838892657669755, 13.89419533942459),
    @(24.58255945752227, 10.35806415387842, 12.96216441575742, 12.59784614067261, 61.1494954677785, 22.01113159693234, 33.57277305830528, 7.839549313819527, 13.86144727313817),
    @(23.87192761984314, 9.617884280056959, 12.96410596608878, 12.6454909462792, 60.83173605650778, 22.05573046941561, 33.68366840751857, 7.840579707167654, 13.83493917711785)
)

$cols = @("B", "C", "D", "E", "G", "H", "I", "J", "L")
$startRow = 2

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $rowValues[$j]
    }
}
